$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (the "Resolving-Mac" sending-cluster rows no longer exist)
$ws.Range("A8:T10").EntireRow.Delete()

# Update remaining rows 2-7: sending cluster "Resolving-Mac" column B values
# shift from Artn(23)->(22) etc. due to removed shared string, and the
# underlying TPM-derived metrics were recomputed.

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.121263666666666
$ws.Range("H2").Value = 3.363791
$ws.Range("I2").Value = 0.8978163344397481
$ws.Range("J2").Value = 0.9294754023256565
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1011536666666667
$ws.Range("N2").Value = 0.303461
$ws.Range("O2").Value = 0.007629860605400263
$ws.Range("P2").Value = 0.008254451482408482
$ws.Range("Q2").Value = 0.1134199311834444
$ws.Range("R2").Value = 1.020779380651
$ws.Range("S2").Value = 0.006850213481026701
$ws.Range("T2").Value = 0.007672309612589236

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.121263666666666
$ws.Range("H3").Value = 3.363791
$ws.Range("I3").Value = 0.8978163344397481
$ws.Range("J3").Value = 0.9294754023256565
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.14695433333333
$ws.Range("N3").Value = 30.440863
$ws.Range("O3").Value = 0.7653686681256785
$ws.Range("P3").Value = 0.8280227993585454
$ws.Range("Q3").Value = 11.37741122129255
$ws.Range("R3").Value = 102.396700991633
$ws.Range("S3").Value = 0.6871604921116287
$ws.Range("T3").Value = 0.7696268245686003

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.121263666666666
$ws.Range("H4").Value = 3.363791
$ws.Range("I4").Value = 0.8978163344397481
$ws.Range("J4").Value = 0.9294754023256565
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.0094955
$ws.Range("N4").Value = 6.018991
$ws.Range("O4").Value = 0.2270014712689213
$ws.Range("P4").Value = 0.1637227491590462
$ws.Range("Q4").Value = 3.374437959146833
$ws.Range("R4").Value = 20.246627754881
$ws.Range("S4").Value = 0.2038056288470927
$ws.Range("T4").Value = 0.152176268144467

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.127615
$ws.Range("H5").Value = 0.25523
$ws.Range("I5").Value = 0.1021836655602519
$ws.Range("J5").Value = 0.07052459767434344
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1011536666666667
$ws.Range("N5").Value = 0.303461
$ws.Range("O5").Value = 0.007629860605400263
$ws.Range("P5").Value = 0.008254451482408482
$ws.Range("Q5").Value = 0.01290872517166667
$ws.Range("R5").Value = 0.07745235103000001
$ws.Range("S5").Value = 0.0007796471243735619
$ws.Range("T5").Value = 0.000582141869819246

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.127615
$ws.Range("H6").Value = 0.25523
$ws.Range("I6").Value = 0.1021836655602519
$ws.Range("J6").Value = 0.07052459767434344
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.14695433333333
$ws.Range("N6").Value = 30.440863
$ws.Range("O6").Value = 0.7653686681256785
$ws.Range("P6").Value = 0.8280227993585454
$ws.Range("Q6").Value = 1.294903577248333
$ws.Range("R6").Value = 7.769421463490001
$ws.Range("S6").Value = 0.07820817601404979
$ws.Range("T6").Value = 0.05839597478994502

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.127615
$ws.Range("H7").Value = 0.25523
$ws.Range("I7").Value = 0.1021836655602519
$ws.Range("J7").Value = 0.07052459767434344
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.0094955
$ws.Range("N7").Value = 6.018991
$ws.Range("O7").Value = 0.2270014712689213
$ws.Range("P7").Value = 0.1637227491590462
$ws.Range("Q7").Value = 0.3840567682325
$ws.Range("R7").Value = 1.53622707293
$ws.Range("S7").Value = 0.02319584242182859
$ws.Range("T7").Value = 0.01154648101457918
